$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update gabSf (I) and gabNum (J) values for rows 2-15
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 9).Value = 0.5   # Column I: gabSf
    $ws.Cells.Item($r, 10).Value = 0.2  # Column J: gabNum
}

# Set the active cell / selection on the sheet view to J2
$ws.Range("J2").Select()
